$d = $word.ActiveDocument
$d.Content.Find.Execute("REALIZACIÓN DEL PROTOCOLO", $true, $false, $false, $false, $false,
                         $true, 1, $false, "REALIZACIÓN DEL CONSENTIMIENTO", 2)
